$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/date range) ---
$ws.Range("A8").Characters(21, 2).Text = "47"
$ws.Range("C9").Characters(27, 10).Text = "11/21/2022"
$ws.Range("C9").Characters(48, 10).Text = "11/27/2022"

# --- Crime statistics table updates (rows 15-29) ---

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 28
$ws.Range("J15").Value = 29
$ws.Range("K15").Value = -3.448275862068
$ws.Range("L15").Value = -15.151515151515
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = -6.666666666666

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 3.571428571428
$ws.Range("I16").Value = 280
$ws.Range("J16").Value = 178
$ws.Range("K16").Value = 57.303370786516
$ws.Range("L16").Value = 57.303370786516
$ws.Range("M16").Value = 6.060606060606
$ws.Range("N16").Value = -64.240102171136

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 22.222222222222
$ws.Range("I17").Value = 392
$ws.Range("J17").Value = 286
$ws.Range("K17").Value = 37.062937062937
$ws.Range("L17").Value = 65.400843881856
$ws.Range("M17").Value = 9.192200557103
$ws.Range("N17").Value = -16.772823779193

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 16
$ws.Range("H18").Value = 45.454545454545
$ws.Range("I18").Value = 149
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 63.736263736263
$ws.Range("L18").Value = 7.194244604316
$ws.Range("M18").Value = -25.125628140703
$ws.Range("N18").Value = -84.543568464730

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -4.166666666666
$ws.Range("I19").Value = 364
$ws.Range("J19").Value = 297
$ws.Range("K19").Value = 22.558922558922
$ws.Range("L19").Value = 33.333333333333
$ws.Range("M19").Value = 47.368421052631
$ws.Range("N19").Value = -10.565110565110

# Row 20
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 180
$ws.Range("I20").Value = 203
$ws.Range("J20").Value = 125
$ws.Range("K20").Value = 62.4
$ws.Range("L20").Value = 123.076923076923
$ws.Range("M20").Value = 181.944444444444
$ws.Range("N20").Value = -49.122807017543

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 131
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = 25.961538461538
$ws.Range("I21").Value = 1421
$ws.Range("J21").Value = 1017
$ws.Range("K21").Value = 39.724680432645
$ws.Range("L21").Value = 48.640167364016
$ws.Range("M21").Value = 21.974248927038
$ws.Range("N21").Value = -54.087237479806

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -75
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -70
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -6.666666666666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 27.272727272727

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = -11.111111111111
$ws.Range("L23").Value = 23.076923076923

# Row 24
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 64
$ws.Range("H24").Value = 4.918032786885
$ws.Range("I24").Value = 728
$ws.Range("J24").Value = 503
$ws.Range("K24").Value = 44.731610337972
$ws.Range("L24").Value = 19.148936170212
$ws.Range("M24").Value = 20.330578512396

# Row 25
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 160
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 29.411764705882
$ws.Range("I25").Value = 499
$ws.Range("J25").Value = 358
$ws.Range("K25").Value = 39.385474860335
$ws.Range("L25").Value = 49.849849849849
$ws.Range("M25").Value = -17.792421746293

# Row 26
$ws.Range("C26").Value = 3
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 51
$ws.Range("J26").Value = 58
$ws.Range("K26").Value = -12.068965517241
$ws.Range("L26").Value = -19.047619047619

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = -38.461538461538
$ws.Range("I27").Value = 91
$ws.Range("J27").Value = 102
$ws.Range("K27").Value = -10.784313725490
$ws.Range("L27").Value = 78.431372549019

# Row 28
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -36.111111111111
$ws.Range("L28").Value = 9.523809523809
$ws.Range("M28").Value = -14.814814814814
$ws.Range("N28").Value = -79.279279279279

# Row 29
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 28
$ws.Range("K29").Value = -32.142857142857
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -20.833333333333
$ws.Range("N29").Value = -79.347826086956
